# Auto-generated Excel COM-interop script to apply cryptos list price/volume update
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "40.118.00"
$ws.Range("E2").Value = "  +0.41%  "

$ws.Range("D3").Value = "2.219.43"
$ws.Range("E3").Value = "  +0.20%  "

$ws.Range("E4").Value = "  +0.06%  "

$cell = $ws.Range("D5")
$cell.NumberFormat = "@"
$cell.Value = "294.68"
$cell.Style = "Normal"
$ws.Range("E5").Value = "  +1.90%  "

$cell = $ws.Range("D6")
$cell.NumberFormat = "@"
$cell.Value = "87.65"
$cell.Style = "Normal"
$ws.Range("E6").Value = "  +0.33%  "

$ws.Range("E7").Value = "  +0.70%  "

$ws.Range("E8").Value = "  +0.09%  "

$cell = $ws.Range("D9")
$cell.NumberFormat = "@"
$cell.Value = "0.472"
$cell.Style = "Normal"
$ws.Range("E9").Value = "  +0.54%  "

$cell = $ws.Range("D10")
$cell.NumberFormat = "@"
$cell.Value = "30.80"
$cell.Style = "Normal"
$ws.Range("E10").Value = "  +1.42%  "

$cell = $ws.Range("D11")
$cell.NumberFormat = "@"
$cell.Value = "51.36"
$cell.Style = "Normal"
$ws.Range("E11").Value = "  +7.55%  "

$cell = $ws.Range("D12")
$cell.NumberFormat = "@"
$cell.Value = "0.0783"
$cell.Style = "Normal"
$ws.Range("E12").Value = "  +1.00%  "

$ws.Range("E13").Value = "  +3.71%  "

$cell = $ws.Range("D14")
$cell.NumberFormat = "@"
$cell.Value = "6.41"
$cell.Style = "Normal"
$ws.Range("E14").Value = "  -0.67%  "

$ws.Range("D15").Value = "2.564.88"
$ws.Range("E15").Value = "  +0.43%  "

$cell = $ws.Range("D16")
$cell.NumberFormat = "@"
$cell.Value = "13.89"
$cell.Style = "Normal"
$ws.Range("E16").Value = "  -0.26%  "

$ws.Range("D17").Value = "2.240.36"
$ws.Range("E17").Value = "  +1.85%  "

$cell = $ws.Range("D18")
$cell.NumberFormat = "@"
$cell.Value = "0.736"
$cell.Style = "Normal"
$ws.Range("E18").Value = "  +1.18%  "

$ws.Range("D19").Value = "40.063.22"
$ws.Range("E19").Value = "  +0.43%  "

$ws.Range("D20").Value = "0.0₃0888"
$ws.Range("E20").Value = "  +0.78%  "

$cell = $ws.Range("D21")
$cell.NumberFormat = "@"
$cell.Value = "11.25"
$cell.Style = "Normal"
$ws.Range("E21").Value = "  -2.04%  "

$ws.Range("E22").Value = "  +0.05%  "

$cell = $ws.Range("D23")
$cell.NumberFormat = "@"
$cell.Value = "65.65"
$cell.Style = "Normal"
$ws.Range("E23").Value = "  +0.17%  "

$cell = $ws.Range("D24")
$cell.NumberFormat = "@"
$cell.Value = "235.59"
$cell.Style = "Normal"
$ws.Range("E24").Value = "  -0.53%  "

$ws.Range("E25").Value = "  -0.01%  "

$ws.Range("E26").Value = "  +1.67%  "

$ws.Range("E27").Value = "  +0.23%  "

$cell = $ws.Range("D28")
$cell.NumberFormat = "@"
$cell.Value = "23.21"
$cell.Style = "Normal"
$ws.Range("E28").Value = "  +3.23%  "

$cell = $ws.Range("D29")
$cell.NumberFormat = "@"
$cell.Value = "9.34"
$cell.Style = "Normal"
$ws.Range("E29").Value = "  +1.58%  "

$ws.Range("E30").Value = "  -4.71%  "

$cell = $ws.Range("D31")
$cell.NumberFormat = "@"
$cell.Value = "161.72"
$cell.Style = "Normal"
$ws.Range("E31").Value = "  +3.77%  "

$cell = $ws.Range("D32")
$cell.NumberFormat = "@"
$cell.Value = "31.75"
$cell.Style = "Normal"
$ws.Range("E32").Value = "  +0.34%  "

$cell = $ws.Range("D33")
$cell.NumberFormat = "@"
$cell.Value = "1.00"
$cell.Style = "Normal"
$ws.Range("E33").Value = "  +0.05%  "

$ws.Range("E34").Value = "  +8.09%  "

$cell = $ws.Range("D35")
$cell.NumberFormat = "@"
$cell.Value = "4.96"
$cell.Style = "Normal"
$ws.Range("E35").Value = "  +0.89%  "

$cell = $ws.Range("D36")
$cell.NumberFormat = "@"
$cell.Value = "0.0715"
$cell.Style = "Normal"
$ws.Range("E36").Value = "  -0.04%  "

$ws.Range("E37").Value = "  -1.44%  "

$ws.Range("E38").Value = "  +1.81%  "

$ws.Range("E39").Value = "  +4.26%  "

$cell = $ws.Range("D40")
$cell.NumberFormat = "@"
$cell.Value = "0.0999"
$cell.Style = "Normal"
$ws.Range("E40").Value = "  +1.35%  "

$cell = $ws.Range("D41")
$cell.NumberFormat = "@"
$cell.Value = "15.64"
$cell.Style = "Normal"
$ws.Range("E41").Value = "  -0.35%  "

$ws.Range("D42").Value = "2.079.74"
$ws.Range("E42").Value = "  -1.18%  "

$cell = $ws.Range("D43")
$cell.NumberFormat = "@"
$cell.Value = "3.75"
$cell.Style = "Normal"
$ws.Range("E43").Value = "  -1.68%  "

$cell = $ws.Range("D44")
$cell.NumberFormat = "@"
$cell.Value = "19.46"
$cell.Style = "Normal"
$ws.Range("E44").Value = "  +11.85%  "

$ws.Range("E45").Value = "  +1.68%  "

$cell = $ws.Range("D46")
$cell.NumberFormat = "@"
$cell.Value = "9.93"
$cell.Style = "Normal"
$ws.Range("E46").Value = "  -0.18%  "

$cell = $ws.Range("D47")
$cell.NumberFormat = "@"
$cell.Value = "2.77"
$cell.Style = "Normal"
$ws.Range("E47").Value = "  +4.83%  "

$ws.Range("E48").Value = "  -10.68%  "

$ws.Range("D49").Value = "2.439.02"
$ws.Range("E49").Value = "  +0.46%  "

$cell = $ws.Range("D50")
$cell.NumberFormat = "@"
$cell.Value = "1.12"
$cell.Style = "Normal"
$ws.Range("E50").Value = "  +3.35%  "

$cell = $ws.Range("D51")
$cell.NumberFormat = "@"
$cell.Value = "1.46"
$cell.Style = "Normal"
$ws.Range("E51").Value = "  +1.61%  "
